# Applies the cryptos.xlsx cell-value diff (prices/volumes refresh + a few row swaps).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value, forcing text storage for pure-numeric-looking strings
# (column D prices like '205.34' or '0.0581') so Excel doesn't silently convert
# them to real numbers -- the source sheet stores every Price/Volume cell as text.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = "26.776.85"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "1.537.50"
$ws.Range("E3").Value = "  -1.76%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue $ws.Range("D5") "205.34"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.77%  "
Set-TextValue $ws.Range("D9") "21.25"
$ws.Range("E9").Value = "  -2.73%  "
Set-TextValue $ws.Range("D10") "0.0581"
$ws.Range("E10").Value = "  -0.54%  "
Set-TextValue $ws.Range("D11") "0.0855"
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("D12").Value = "1.756.27"
$ws.Range("E12").Value = "  -1.78%  "
$ws.Range("D13").Value = "1.537.82"
$ws.Range("E13").Value = "  -1.77%  "
$ws.Range("E14").Value = "  -1.53%  "
Set-TextValue $ws.Range("D15") "0.508"
$ws.Range("E15").Value = "  -1.32%  "
$ws.Range("D16").Value = "26.772.96"
$ws.Range("E16").Value = "  -0.05%  "
Set-TextValue $ws.Range("D17") "60.91"
$ws.Range("E17").Value = "  -1.04%  "
Set-TextValue $ws.Range("D18") "212.98"
$ws.Range("E18").Value = "  -0.75%  "
Set-TextValue $ws.Range("D19") "7.23"
$ws.Range("E19").Value = "  -1.65%  "
$ws.Range("E20").Value = "  +0.61%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("E22").Value = "  -1.82%  "
Set-TextValue $ws.Range("D23") "9.13"
$ws.Range("E23").Value = "  -1.98%  "
Set-TextValue $ws.Range("D24") "1.94"
$ws.Range("E24").Value = "  -3.10%  "
Set-TextValue $ws.Range("D25") "151.70"
$ws.Range("E25").Value = "  -0.70%  "
$ws.Range("E26").Value = "  -2.50%  "
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  -0.98%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D30") "0.0457"
$ws.Range("E30").Value = "  -1.30%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D31") "1.10"
$ws.Range("E31").Value = "  -1.52%  "
$ws.Range("E32").Value = "  +2.29%  "
$ws.Range("D33").Value = "1.363.84"
$ws.Range("E33").Value = "  -1.58%  "
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("E35").Value = "  -2.72%  "
Set-TextValue $ws.Range("D36") "0.964"
$ws.Range("E36").Value = "  +4.50%  "
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("E38").Value = "  +0.53%  "
$ws.Range("E39").Value = "  -1.21%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D40") "5.74"
$ws.Range("E40").Value = "  +7.51%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D41") "0.802"
$ws.Range("E41").Value = "  -1.90%  "
$ws.Range("E42").Value = "  +0.36%  "
Set-TextValue $ws.Range("D43") "2.20"
$ws.Range("E43").Value = "  +0.75%  "
Set-TextValue $ws.Range("D44") "62.92"
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("E45").Value = "  -3.28%  "
$ws.Range("B46").Value = "mCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
Set-TextValue $ws.Range("D46") "2.25"
$ws.Range("E46").Value = "  -3.95%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.670.76"
$ws.Range("E47").Value = "  -1.73%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D48") "85.23"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D49") "0.0508"
$ws.Range("E49").Value = "  +3.29%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₇0976"
$ws.Range("E50").Value = "  -1.38%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D51") "0.0943"
$ws.Range("E51").Value = "  -0.94%  "
